$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.479.71"
$ws.Range("E2").Value = "  +0.82%  "

# Row 3
$ws.Range("D3").Value = "2.512.88"
$ws.Range("E3").Value = "  +2.07%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.96"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.57"
$ws.Range("E6").Value = "  -0.40%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.38%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("E8").Value = "  +0.37%  "

# Row 9
$ws.Range("D9").Value = "2.510.67"
$ws.Range("E9").Value = "  +1.66%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0975"
$ws.Range("E10").Value = "  -0.60%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  -0.70%  "

# Row 12
$ws.Range("E12").Value = "  -2.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.332"
$ws.Range("E13").Value = "  -1.96%  "

# Row 14
$ws.Range("D14").Value = "2.958.97"
$ws.Range("E14").Value = "  +2.08%  "

# Row 15
$ws.Range("D15").Value = "58.419.83"
$ws.Range("E15").Value = "  +0.82%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.08"
$ws.Range("E16").Value = "  -0.88%  "

# Row 17
$ws.Range("E17").Value = "  -0.35%  "

# Row 18
$ws.Range("D18").Value = "2.510.72"
$ws.Range("E18").Value = "  +1.61%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  -0.10%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.30"
$ws.Range("E20").Value = "  +0.26%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$ws.Range("E22").Value = "  +7.66%  "

# Row 23
$ws.Range("E23").Value = "  +0.10%  "

# Row 24
$ws.Range("E24").Value = "  -0.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.403"
$ws.Range("E25").Value = "  -1.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  +0.44%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("E28").Value = "  +0.32%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0753"
$ws.Range("E29").Value = "  +0.80%  "

# Row 30
$ws.Range("E30").Value = "  +1.82%  "

# Row 31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +2.75%  "

# Row 32
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.82"
$ws.Range("E32").Value = "  +0.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.25"
$ws.Range("E33").Value = "  +0.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.06"
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.24"
$ws.Range("E37").Value = "  -7.77%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("E38").Value = "  -0.28%  "

# Row 39
$ws.Range("E39").Value = "  +0.01%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.16"
$ws.Range("E40").Value = "  -0.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.775"
$ws.Range("E41").Value = "  -2.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "277.98"
$ws.Range("E42").Value = "  +1.29%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.48"
$ws.Range("E43").Value = "  +0.86%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  +0.74%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.595"
$ws.Range("E45").Value = "  +1.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "122.86"
$ws.Range("E46").Value = "  -0.91%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0919"
$ws.Range("E47").Value = "  +1.18%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0500"
$ws.Range("E48").Value = "  +2.48%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.66"
$ws.Range("E49").Value = "  +0.36%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0213"
$ws.Range("E50").Value = "  +0.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.83"
$ws.Range("E51").Value = "  -0.24%  "
